$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural edits -------------------------------------------------
# Remove the "Adresse pile relative a ESP" column (old column E) entirely.
$ws.Columns("E").Delete()

# Remove the last row (old row 14, "EBP-0x20") - table now ends at row 13.
$ws.Range("A14").EntireRow.Delete()

# --- Content edits ------------------------------------------------------
# Fix the casing of the absolute stack address next to "Saved EBP".
$ws.Range("D5").Value = "0xffffd65c"

# Fill in the newly-visible "Adresse pile absolue" column for the local
# variables block (rows 7-13) with the real stack addresses.
$ws.Range("D7").Value = "0xffffd654"
$ws.Range("D8").Value = "0xffffd650"
$ws.Range("D9").Value = "0xffffd64c"
$ws.Range("D10").Value = "0xffffd648"
$ws.Range("D11").Value = "0xffffd644"
$ws.Range("D12").Value = "0xffffd640"
$ws.Range("D13").Value = "0xffffd63c"

# --- Formatting -----------------------------------------------------
# Every cell in the table gets the bordered / centred look (the merged
# B7:B13 cell keeps its own centred+middle style untouched).
$ws.Range("A4:D13").HorizontalAlignment = -4108

# Column width tweaks for columns A and B.
$ws.Columns("A").ColumnWidth = 16.5
$ws.Columns("B").ColumnWidth = 18

# --- View -------------------------------------------------------------
$ws.Range("A4:D13").Select()
